$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F40").Value2 = 117.000
$ws.Range("H40").Value2 = 5496.05
$ws.Range("H42").Value2 = 6914.88
$ws.Range("F44").Value2 = 49.000
$ws.Range("H44").Value2 = 27379.98
$ws.Range("F54").Value2 = 716.000
$ws.Range("H54").Value2 = 6512.02
$ws.Range("F56").Value2 = 2707.000
$ws.Range("H56").Value2 = 27832.83
$ws.Range("H68").Value2 = 6017.23
$ws.Range("F80").Value2 = 1255.000
$ws.Range("H80").Value2 = 7714.01
$ws.Range("F82").Value2 = 2646.000
$ws.Range("H82").Value2 = 68922.04
$ws.Range("F96").Value2 = 248.000
$ws.Range("H96").Value2 = 20993.57
$ws.Range("F97").Value2 = 2082.000
$ws.Range("H97").Value2 = 286785.09
$ws.Range("F113").Value2 = 596.000
$ws.Range("H113").Value2 = 22118.06
$ws.Range("F115").Value2 = 1051.000
$ws.Range("H115").Value2 = 34627.03
$ws.Range("F120").Value2 = 12877.000
$ws.Range("H120").Value2 = 49984.49
$ws.Range("F174").Value2 = 2454.000
$ws.Range("H174").Value2 = 45951.15
$ws.Range("F175").Value2 = 109836.000
$ws.Range("H175").Value2 = 166154.32
$ws.Range("F178").Value2 = 775.000
$ws.Range("H178").Value2 = 135328.15
$ws.Range("F182").Value2 = 126.000
$ws.Range("H182").Value2 = 58447.08
$ws.Range("H189").Value2 = 6893.93
$ws.Range("F231").Value2 = 20012.800
$ws.Range("H231").Value2 = 45414.96
$ws.Range("F235").Value2 = 325.600
$ws.Range("H235").Value2 = 3249.21
$ws.Range("H245").Value2 = 1188.56
$ws.Range("F247").Value2 = 5364.000
$ws.Range("H247").Value2 = 13465.57
